$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.729.24'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.885.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4759'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.87%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2846'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06545'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +15.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.884.47'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07560'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '95.51'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +13.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.074'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6500'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '303.94'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +32.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.723.44'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9993'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007551'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.127.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9987'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.133'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.155'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '168.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.227'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +13.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.947'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1075'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.349'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.154'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.86%  '
$ws.Range('E32').Value = '  +3.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05048'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.170'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7232'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.716'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01924'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.710'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.071'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8993'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.72%  '
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9995'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4197'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.616'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.68%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.88'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.41%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.338'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.52%  '
$ws.Range('E47').Value = '  +1.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.960'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05599'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.383'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.65%  '
